# Update the player roster table on Sheet1 (A2:C19) to reflect the new
# lineup of players, positions and teams.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Derrick White"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Boston Celtics"

$ws.Range("A3").Value = "Coby White"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Chicago Bulls"

$ws.Range("A4").Value = "Cade Cunningham"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Detroit Pistons"

$ws.Range("A5").Value = "Shaedon Sharpe"
$ws.Range("B5").Value = "SG,SF"
$ws.Range("C5").Value = "Portland Trail Blazers"

$ws.Range("A6").Value = "Ausar Thompson"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Detroit Pistons"

$ws.Range("A7").Value = "Naz Reid"
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "Minnesota Timberwolves"

$ws.Range("A8").Value = "Aaron Nesmith"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Indiana Pacers"

$ws.Range("A9").Value = "Onyeka Okongwu"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Atlanta Hawks"

$ws.Range("A10").Value = "Karlo Matkovic"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "New Orleans Pelicans"

$ws.Range("A11").Value = "Isaiah Hartenstein"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Oklahoma City Thunder"

$ws.Range("A12").Value = "Damian Lillard"
$ws.Range("B12").Value = "PG"
$ws.Range("C12").Value = "Milwaukee Bucks"

$ws.Range("A13").Value = "Jusuf Nurkic"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Charlotte Hornets"

$ws.Range("A14").Value = "Gary Trent Jr."
$ws.Range("B14").Value = "PG,SG,SF"
$ws.Range("C14").Value = "Milwaukee Bucks"

$ws.Range("A15").Value = "LaMelo Ball"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Charlotte Hornets"

$ws.Range("A16").Value = "Malik Monk"
$ws.Range("B16").Value = "PG,SG,SF"
$ws.Range("C16").Value = "Sacramento Kings"

$ws.Range("A17").Value = "Collin Sexton"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("C17").Value = "Utah Jazz"

$ws.Range("A18").Value = "Anthony Davis"
$ws.Range("B18").Value = "PF,C"
$ws.Range("C18").Value = "Dallas Mavericks"

$ws.Range("A19").Value = "Andrew Wiggins"
$ws.Range("B19").Value = "SF,PF"
$ws.Range("C19").Value = "Miami Heat"
